$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.708.59'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '1.868.41'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7286'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3122'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07083'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.36'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08230'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7446'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.305'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").Value = '1.858.29'
$ws.Range("E14").Value = '  -0.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.32%  '

$ws.Range("D16").Value = '29.708.47'
$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.998'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '247.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007786'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.13%  '

$ws.Range("D22").Value = '2.110.90'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.702'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.58%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1528'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.157'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.009'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.435'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.508'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.40%  '

$ws.Range("E32").Value = '  -0.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.182'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05266'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.229'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7536'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9977'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.694'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01929'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.736'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4458'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.986'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8660'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.93'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.74%  '

$ws.Range("D45").Value = '1.046.66'
$ws.Range("E45").Value = '  -5.46%  '

$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.464'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.815'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.49%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.006.67'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.876'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.77%  '

